# "Adicionar e Vender" — register new purchases of concrete products
# (Pedigree / Coelho) replacing the generic placeholder row, update the
# Compras / Estoque ledgers accordingly, and drop the unused "Unidade"
# method.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 - "Métodos": remove the unused "Unidade" method (row 3)
# ---------------------------------------------------------------------
$wsMetodos = $wb.Worksheets.Item(1)
$wsMetodos.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# Sheet 2 - "Produtos": replace the generic "Cachorro" product with two
# concrete products, Pedigree and Coelho
# ---------------------------------------------------------------------
$wsProdutos = $wb.Worksheets.Item(2)
$wsProdutos.Cells.Item(2,2).Value2 = "Pedigree"
$wsProdutos.Cells.Item(2,5).Value2 = 70
$wsProdutos.Cells.Item(2,6).Value2 = 80

$wsProdutos.Cells.Item(3,1).Value2 = "Ração"
$wsProdutos.Cells.Item(3,2).Value2 = "Coelho"
$wsProdutos.Cells.Item(3,3).Value2 = "Pacote"
$wsProdutos.Cells.Item(3,4).Value2 = "Pacote"
$wsProdutos.Cells.Item(3,5).Value2 = 60
$wsProdutos.Cells.Item(3,6).Value2 = 90
$wsProdutos.Cells.Item(3,7).Value2 = "Não"

# ---------------------------------------------------------------------
# Sheet 5 - "P_Compras": new purchase detail lines
# ---------------------------------------------------------------------
$wsPCompras = $wb.Worksheets.Item(5)

$wsPCompras.Cells.Item(2,1).Value2 = 10000
$wsPCompras.Cells.Item(2,2).Value2 = "Ração"
$wsPCompras.Cells.Item(2,3).Value2 = "Pedigree"
$wsPCompras.Cells.Item(2,4).Value2 = "Pacote"
$wsPCompras.Cells.Item(2,5).Value2 = 20
$wsPCompras.Cells.Item(2,6).Value2 = 70
$wsPCompras.Cells.Item(2,7).Value2 = 1400

$wsPCompras.Cells.Item(3,1).Value2 = 10000
$wsPCompras.Cells.Item(3,2).Value2 = "Ração"
$wsPCompras.Cells.Item(3,3).Value2 = "Coelho"
$wsPCompras.Cells.Item(3,4).Value2 = "Pacote"
$wsPCompras.Cells.Item(3,5).Value2 = 10
$wsPCompras.Cells.Item(3,6).Value2 = 60
$wsPCompras.Cells.Item(3,7).Value2 = 600

# Row 4 is written as text ("quote-prefixed" style of entry) - values
# are kept as literal text rather than numbers
$rngRow4 = $wsPCompras.Range("A4:G4")
$rngRow4.NumberFormat = "@"
$wsPCompras.Cells.Item(4,1).Value = "10001"
$wsPCompras.Cells.Item(4,2).Value = "Ração"
$wsPCompras.Cells.Item(4,3).Value = "Pedigree"
$wsPCompras.Cells.Item(4,4).Value = "Pacote"
$wsPCompras.Cells.Item(4,5).Value = "10"
$wsPCompras.Cells.Item(4,6).Value = "70.0"
$wsPCompras.Cells.Item(4,7).Value = "700.0"

# ---------------------------------------------------------------------
# Sheet 6 - "Compras": purchase order header lines
# ---------------------------------------------------------------------
$wsCompras = $wb.Worksheets.Item(6)

$wsCompras.Cells.Item(2,1).Value2 = 10000
$wsCompras.Cells.Item(2,2).Value2 = 0
$wsCompras.Cells.Item(2,3).Value2 = "17/01/2023"
$wsCompras.Cells.Item(2,4).Value2 = 0

$wsCompras.Cells.Item(3,1).Value2 = 10000
$wsCompras.Cells.Item(3,2).Value2 = 30
$wsCompras.Cells.Item(3,3).Value2 = "17/01/2023"
$wsCompras.Cells.Item(3,4).Value2 = 2000

$wsCompras.Cells.Item(4,1).Value2 = 10001
$wsCompras.Cells.Item(4,2).Value2 = 0
$wsCompras.Cells.Item(4,3).Value2 = "17/01/2023"
$wsCompras.Cells.Item(4,4).Value2 = 0

# ---------------------------------------------------------------------
# Sheet 7 - "Estoque": stock now split across the two concrete products
# ---------------------------------------------------------------------
$wsEstoque = $wb.Worksheets.Item(7)

$wsEstoque.Cells.Item(2,1).Value2 = "Ração"
$wsEstoque.Cells.Item(2,2).Value2 = "Pedigree"
$wsEstoque.Cells.Item(2,3).Value2 = "Pacote"
$wsEstoque.Cells.Item(2,4).Value2 = 30

$wsEstoque.Cells.Item(3,1).Value2 = "Ração"
$wsEstoque.Cells.Item(3,2).Value2 = "Coelho"
$wsEstoque.Cells.Item(3,3).Value2 = "Pacote"
$wsEstoque.Cells.Item(3,4).Value2 = 10
